# "Generate Report for Handback"
# This run stamps fresh handoff/handback timestamps for the
# bb7ca973-2619-41a2-a4d2-ed0759aca5bf file (row 3 in each language sheet),
# and (for de-de) refreshes the handoff timestamp that happens to equal the
# shared "Latest HO Xliff Generate Date" on the Overview sheet.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: bump the "Latest HO Xliff Generate Date" for the
#     bb7ca973... row (row 3) to reflect the newly generated handback report.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-12 16:59:51"

# --- zh-cn sheet: row 3 is the bb7ca973... file. Its handoff file was
#     regenerated, so both the "Correspond Handoff Datetime" (H) and the
#     "Correspond Handback DateTime" (K) move forward.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-12 16:59:44"
$wsZhCn.Range("K3").Value = "2016-08-12 17:00:30"

# --- de-de sheet: row 3 is the bb7ca973... file; same refresh of the
#     handoff/handback timestamps.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-12 16:59:51"
$wsDeDe.Range("K3").Value = "2016-08-12 17:00:40"
